# The sheet holds monthly index values in rows 2..49 (A=yyyy-mm label,
# B/C = numeric index values), laid out in 12-row (one per year) blocks
# that originally run January..December. The edit re-orders every
# 12-row year-block so it instead runs October, November, December,
# January, February, ... September (i.e. Oct-Dec are rotated to the
# front of each year's block). Header row 1 and everything else is
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 49
$yearBlockSize = 12

# 1) Snapshot the current (pre-edit) rows into memory so writes below
#    don't clobber values we still need to read. (Value2 is used rather
#    than Value to get the raw primitive without variant/date wrapping.)
$rows = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $rows[$r] = @{
        A = $ws.Cells.Item($r, 1).Value2
        B = $ws.Cells.Item($r, 2).Value2
        C = $ws.Cells.Item($r, 3).Value2
    }
}

# 2) For each 12-row year block, rotate so the last 3 rows (Oct, Nov,
#    Dec) move to the front, followed by the first 9 rows (Jan..Sep),
#    preserving their relative order.
$blockStart = $firstDataRow
while ($blockStart -le $lastDataRow) {
    $blockEnd = [Math]::Min($blockStart + $yearBlockSize - 1, $lastDataRow)

    $sourceRows = @($blockStart..$blockEnd)
    $rotateCount = [Math]::Min(3, $sourceRows.Count)
    $tail = $sourceRows[($sourceRows.Count - $rotateCount)..($sourceRows.Count - 1)]
    $head = $sourceRows[0..($sourceRows.Count - $rotateCount - 1)]
    $newOrder = $tail + $head

    $destRow = $blockStart
    foreach ($srcRow in $newOrder) {
        $ws.Range("A$destRow").Value2 = $rows[$srcRow].A
        $ws.Range("B$destRow").Value2 = $rows[$srcRow].B
        $ws.Range("C$destRow").Value2 = $rows[$srcRow].C
        $destRow++
    }

    $blockStart += $yearBlockSize
}
